# Update the figure in T2 on the active sheet (row 2), matching the
# refreshed daily export value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("T2").Value = 502242
